$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 616.34784
$ws.Cells.Item(28, 9).Value = 391.1111
$ws.Cells.Item(28, 10).Value = 1427.2
$ws.Cells.Item(28, 11).Value = 391.1111
$ws.Cells.Item(28, 12).Value = 1427.2
$ws.Cells.Item(28, 13).Value = 93.88889999999998
$ws.Cells.Item(28, 14).Value = -2397.2
$ws.Cells.Item(40, 8).Value = 4731.3335
$ws.Cells.Item(40, 9).Value = 2997.5
$ws.Cells.Item(40, 10).Value = 4998.077
$ws.Cells.Item(40, 11).Value = 2997.5
$ws.Cells.Item(40, 12).Value = 4998.077
$ws.Cells.Item(40, 13).Value = -2822.5
$ws.Cells.Item(40, 14).Value = -5348.077
$ws.Cells.Item(62, 8).Value = 6340.231
$ws.Cells.Item(62, 9).Value = 4697
$ws.Cells.Item(62, 11).Value = 4697
$ws.Cells.Item(62, 13).Value = -4073
$ws.Cells.Item(65, 8).Value = 6340.231
$ws.Cells.Item(65, 9).Value = 4697
$ws.Cells.Item(65, 11).Value = 23485
$ws.Cells.Item(65, 13).Value = -20365
$ws.Cells.Item(76, 8).Value = 4997
$ws.Cells.Item(76, 9).Value = 5003
$ws.Cells.Item(76, 11).Value = 5003
$ws.Cells.Item(76, 13).Value = -4688
$ws.Cells.Item(79, 8).Value = 4997
$ws.Cells.Item(79, 9).Value = 5003
$ws.Cells.Item(79, 11).Value = 5003
$ws.Cells.Item(79, 13).Value = -3911
$ws.Cells.Item(116, 8).Value = 2807
$ws.Cells.Item(116, 9).Value = 2774.25
$ws.Cells.Item(116, 10).Value = 2872.5
$ws.Cells.Item(116, 11).Value = 2774.25
$ws.Cells.Item(116, 12).Value = 2872.5
$ws.Cells.Item(116, 13).Value = 667.75
$ws.Cells.Item(116, 14).Value = -9756.5
$ws.Cells.Item(124, 8).Value = 104228
$ws.Cells.Item(124, 10).Value = 104228
$ws.Cells.Item(124, 12).Value = 104228
$ws.Cells.Item(124, 14).Value = -114048
$ws.Cells.Item(126, 8).Value = 94997.836
$ws.Cells.Item(126, 10).Value = 94997.836
$ws.Cells.Item(126, 12).Value = 94997.836
$ws.Cells.Item(126, 14).Value = -104877.836
$ws.Cells.Item(138, 8).Value = 269570.84
$ws.Cells.Item(138, 9).Value = 31105.908
$ws.Cells.Item(138, 10).Value = 1253238.8
$ws.Cells.Item(138, 11).Value = 93317.724
$ws.Cells.Item(138, 12).Value = 3759716.4
$ws.Cells.Item(138, 13).Value = -88177.724
$ws.Cells.Item(138, 14).Value = -3769996.4
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5462.8525
$ws.Cells.Item(32, 9).Value = 4911.7964
$ws.Cells.Item(32, 11).Value = 4911.7964
$ws.Cells.Item(32, 13).Value = -4624.7964
$ws.Cells.Item(45, 8).Value = 2691.2222
$ws.Cells.Item(45, 9).Value = 2540.6365
$ws.Cells.Item(45, 10).Value = 2927.8572
$ws.Cells.Item(45, 11).Value = 2540.6365
$ws.Cells.Item(45, 12).Value = 2927.8572
$ws.Cells.Item(45, 13).Value = -2163.6365
$ws.Cells.Item(45, 14).Value = -3681.8572
$ws.Cells.Item(88, 8).Value = 1533.3572
$ws.Cells.Item(88, 9).Value = 1209.375
$ws.Cells.Item(88, 11).Value = 1209.375
$ws.Cells.Item(88, 13).Value = -803.375
$ws.Cells.Item(91, 8).Value = 1533.3572
$ws.Cells.Item(91, 9).Value = 1209.375
$ws.Cells.Item(91, 11).Value = 1209.375
$ws.Cells.Item(91, 13).Value = 194.625
$ws.Cells.Item(122, 8).Value = 1905.2354
$ws.Cells.Item(122, 9).Value = 1774.375
$ws.Cells.Item(122, 11).Value = 5323.125
$ws.Cells.Item(122, 13).Value = -2873.125
$ws.Cells.Item(132, 8).Value = 1713.2258
$ws.Cells.Item(132, 9).Value = 1196.7858
$ws.Cells.Item(132, 11).Value = 3590.3574
$ws.Cells.Item(132, 13).Value = -1060.3574
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 2217.147
$ws.Cells.Item(20, 9).Value = 1818.7727
$ws.Cells.Item(20, 10).Value = 2947.5
$ws.Cells.Item(20, 11).Value = 1818.7727
$ws.Cells.Item(20, 12).Value = 2947.5
$ws.Cells.Item(20, 13).Value = -1571.7727
$ws.Cells.Item(20, 14).Value = -3441.5
$ws.Cells.Item(52, 8).Value = 62497
$ws.Cells.Item(52, 10).Value = 62497
$ws.Cells.Item(52, 12).Value = 62497
$ws.Cells.Item(52, 14).Value = -63023
$ws.Cells.Item(86, 8).Value = 2392
$ws.Cells.Item(86, 9).Value = 2569.6365
$ws.Cells.Item(86, 11).Value = 2569.6365
$ws.Cells.Item(86, 13).Value = -1446.6365
$ws.Cells.Item(89, 8).Value = 2392
$ws.Cells.Item(89, 9).Value = 2569.6365
$ws.Cells.Item(89, 11).Value = 12848.1825
$ws.Cells.Item(89, 13).Value = -7232.182500000001
$ws.Cells.Item(121, 8).Value = 62497
$ws.Cells.Item(121, 10).Value = 62497
$ws.Cells.Item(121, 12).Value = 62497
$ws.Cells.Item(121, 14).Value = -65991
$ws.Cells.Item(134, 8).Value = 6237.357
$ws.Cells.Item(134, 9).Value = 2490
$ws.Cells.Item(134, 11).Value = 7470
$ws.Cells.Item(134, 13).Value = -4935
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2951.5557
$ws.Cells.Item(31, 9).Value = 1517.875
$ws.Cells.Item(31, 11).Value = 1517.875
$ws.Cells.Item(31, 13).Value = -1222.875
$ws.Cells.Item(34, 8).Value = 2951.5557
$ws.Cells.Item(34, 9).Value = 1517.875
$ws.Cells.Item(34, 11).Value = 1517.875
$ws.Cells.Item(34, 13).Value = -1315.875
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(17, 8).Value = 194.625
$ws.Cells.Item(17, 9).Value = 224.5
$ws.Cells.Item(17, 10).Value = 105
$ws.Cells.Item(17, 11).Value = 673.5
$ws.Cells.Item(17, 12).Value = 315
$ws.Cells.Item(17, 13).Value = -504.5
$ws.Cells.Item(17, 14).Value = -653
$ws.Cells.Item(114, 8).Value = 27779772
$ws.Cells.Item(114, 9).Value = 62500170
$ws.Cells.Item(114, 10).Value = 3457
$ws.Cells.Item(114, 11).Value = 187500510
$ws.Cells.Item(114, 12).Value = 10371
$ws.Cells.Item(114, 13).Value = -187497256
$ws.Cells.Item(114, 14).Value = -16879
$ws.Cells.Item(137, 8).Value = 2005.1052
$ws.Cells.Item(137, 10).Value = 3199.889
$ws.Cells.Item(137, 12).Value = 9599.667000000001
$ws.Cells.Item(137, 14).Value = -19799.667
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(43, 8).Value = 58576.6
$ws.Cells.Item(43, 9).Value = 31348.8
$ws.Cells.Item(43, 11).Value = 31348.8
$ws.Cells.Item(43, 13).Value = -31197.8
$ws.Cells.Item(113, 8).Value = 2592.182
$ws.Cells.Item(113, 9).Value = 2312.625
$ws.Cells.Item(113, 10).Value = 3337.6667
$ws.Cells.Item(113, 11).Value = 2312.625
$ws.Cells.Item(113, 12).Value = 3337.6667
$ws.Cells.Item(113, 13).Value = -142.625
$ws.Cells.Item(113, 14).Value = -7677.6667
$ws.Cells.Item(132, 8).Value = 9527349
$ws.Cells.Item(132, 9).Value = 16670999
$ws.Cells.Item(132, 11).Value = 50012997
$ws.Cells.Item(132, 13).Value = -50010467
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 5526.2188
$ws.Cells.Item(46, 9).Value = 2537.25
$ws.Cells.Item(46, 11).Value = 2537.25
$ws.Cells.Item(46, 13).Value = -2349.25
$ws.Cells.Item(122, 8).Value = 4371.5713
$ws.Cells.Item(122, 9).Value = 4208.5
$ws.Cells.Item(122, 10).Value = 4519.8184
$ws.Cells.Item(122, 11).Value = 12625.5
$ws.Cells.Item(122, 12).Value = 13559.4552
$ws.Cells.Item(122, 13).Value = -10175.5
$ws.Cells.Item(122, 14).Value = -18459.4552
$ws.Cells.Item(128, 8).Value = 81666.336
$ws.Cells.Item(128, 10).Value = 81666.336
$ws.Cells.Item(128, 12).Value = 81666.336
$ws.Cells.Item(128, 14).Value = -91626.336
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 47354.72
$ws.Cells.Item(81, 9).Value = 116361.11
$ws.Cells.Item(81, 10).Value = 8538.625
$ws.Cells.Item(81, 11).Value = 232722.22
$ws.Cells.Item(81, 12).Value = 17077.25
$ws.Cells.Item(81, 13).Value = -231661.22
$ws.Cells.Item(81, 14).Value = -19199.25
$ws.Cells.Item(84, 8).Value = 47354.72
$ws.Cells.Item(84, 9).Value = 116361.11
$ws.Cells.Item(84, 10).Value = 8538.625
$ws.Cells.Item(84, 11).Value = 1163611.1
$ws.Cells.Item(84, 12).Value = 85386.25
$ws.Cells.Item(84, 13).Value = -1158307.1
$ws.Cells.Item(84, 14).Value = -95994.25
$ws.Cells.Item(121, 8).Value = 164113.8
$ws.Cells.Item(121, 10).Value = 164113.8
$ws.Cells.Item(121, 12).Value = 164113.8
$ws.Cells.Item(121, 14).Value = -167607.8
$ws.Cells.Item(122, 8).Value = 3101.2917
$ws.Cells.Item(122, 9).Value = 2896.738
$ws.Cells.Item(122, 11).Value = 8690.214
$ws.Cells.Item(122, 13).Value = -6240.214
